# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) updates -----------------------------
# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.619.60"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.643.09"
$ws.Range("E3").Value = "  +0.67%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.19%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'215.83"
$ws.Range("E5").Value = "  +1.38%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.86%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.16%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.07%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.77%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'19.21"
$ws.Range("E10").Value = "  +0.44%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.15%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.873.12"
$ws.Range("E12").Value = "  +0.69%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +3.43%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.642.20"
$ws.Range("E14").Value = "  +0.95%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.531"
$ws.Range("E15").Value = "  +1.52%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'65.86"
$ws.Range("E16").Value = "  +4.30%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.668.08"
$ws.Range("E17").Value = "  +0.11%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +1.51%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'218.16"
$ws.Range("E19").Value = "  +0.26%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.28%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +2.15%  "

# Row 22 - Chainlink
$ws.Range("E22").Value = "  +2.07%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "'9.55"
$ws.Range("E23").Value = "  +2.03%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +11.39%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'146.33"
$ws.Range("E25").Value = "  -1.10%  "

# Row 26 - BinanceUSD (unchanged)

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.17%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  +3.64%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +2.60%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +2.77%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.07%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'3.39"
$ws.Range("E32").Value = "  +3.17%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +2.73%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.276.86"
$ws.Range("E34").Value = "  +5.64%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +2.54%  "

# Row 36 - VeChain
$ws.Range("E36").Value = "  +6.24%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +0.24%  "

# Rows 38/39 - ARBITRUM and ImmutableX swap places (ImmutableX now ranks
# above ARBITRUM), each with refreshed price/volume figures.
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.527"
$ws.Range("E38").Value = "  +5.34%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'0.827"
$ws.Range("E39").Value = "  +2.51%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  +0.22%  "

# Row 41 - TrustWalletToken
$ws.Range("E41").Value = "  +2.05%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  -1.36%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  +0.94%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.784.07"
$ws.Range("E44").Value = "  +0.64%  "

# Row 45 - Quant
$ws.Range("D45").Value = "'93.08"
$ws.Range("E45").Value = "  +0.69%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'59.72"
$ws.Range("E46").Value = "  +9.39%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  +3.58%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +0.62%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "'7.79"
$ws.Range("E49").Value = "  +2.40%  "

# Row 50 - Algorand
$ws.Range("D50").Value = "'0.0976"
$ws.Range("E50").Value = "  +4.02%  "

# Row 51 - Mantle
$ws.Range("E51").Value = "  -0.69%  "
